# Add new column 'Correction' to Card18
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card18")

# Header cell N1: new "Correction" column header, formatted like the rest
# of the header row (bold font + border -> copy M1's format via
# PasteSpecial, since a plain `.Style =` assignment only resets formatting
# to the default style rather than copying it).
$ws.Cells.Item(1, 13).Copy()
$ws.Cells.Item(1, 14).PasteSpecial(-4122)
$ws.Cells.Item(1, 14).Value = "Correction"
$excel.CutCopyMode = $false

# Existing column M ("event") cells were blank inline strings; they now
# carry literal "nan" text like the other empty-valued columns.
# Column N ("Correction") is added as blank (empty string) cells for every
# data row, matching the plain (unstyled) data-row formatting.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 13).Value = "nan"

    # A bare "" assignment deletes the cell instead of leaving an empty
    # text cell behind, so force text entry with a leading apostrophe and
    # then reset the resulting quote-prefix formatting back to the plain
    # (default) style used by the rest of the data rows.
    $ws.Cells.Item($r, 14).Value = "'"
    $ws.Cells.Item($r, 14).Style = $ws.Cells.Item($r, 14).Style
}
